{"js": "// Change the salutation from \"To whom it may concern,\" to \"dear hiring manager,\"\n// and change \"Are you looking for a [job title] with:\" to\n// \"Are you looking for a/an [job title] with:\".\n\nconst body = context.document.body;\n\n// 1) Salutation heading.\nconst heading = body.search(\"To whom it may concern,\", { matchCase: true });\nheading.load(\"text\");\nawait context.sync();\n\nif (heading.items.length > 0) {\n  heading.items[0].insertText(\"dear hiring manager,\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Are you looking for a [job title] with:\" -> add \"/an\" after \"a\".\nconst intro = body.search(\"Are you looking for a [job title] with:\", { matchCase: true });\nintro.load(\"text\");\nawait context.sync();\n\nif (intro.items.length > 0) {\n  intro.items[0].insertText(\"Are you looking for a/an [job title] with:\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Change the salutation from \"To whom it may concern,\" to \"dear hiring manager,\"\n# and change \"Are you looking for a [job title] with:\" to\n# \"Are you looking for a/an [job title] with:\".\n\n$d = $word.ActiveDocument\n\n# 1) Salutation heading.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"To whom it may concern,\"\n$find1.Replacement.Text = \"dear hiring manager,\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# 2) \"Are you looking for a [job title] with:\" -> add \"/an\" after \"a\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Are you looking for a [job title] with:\"\n$find2.Replacement.Text = \"Are you looking for a/an [job title] with:\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
